# Upload new version with timestamp
# -----------------------------------------------------------------------
# A new low-stock item ("HEMOCLAR 0.5% CREAM 40 GM") needs to be inserted
# into the report table (alphabetically, between "DICLAC ..." and
# "KADEE ..."), which lives on row 11. Every row from 11 downwards
# (existing items, the subtotal row and the footer row) needs to shift
# down by one row, the subtotal needs to grow by the new item's price,
# and the footer timestamp needs to be refreshed.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Shift rows 16..11 down to 17..12 (processing bottom-up so that we
#    never overwrite a source row before it has been copied).
#    For every row we: copy the formatting (so styles/number-formats are
#    reused exactly as-is, matching the row being vacated) and then copy
#    over the values of the cells that actually hold data.
#    NOTE: uses positional parameters only - named parameters on custom
#    functions are unreliable in this interpreter.
# ------------------------------------------------------------------

function Move-RowDown($SrcRow, $DstRow, $Cols) {
    $srcRange = "A$SrcRow`:Q$SrcRow"
    $dstRange = "A$DstRow`:Q$DstRow"
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    foreach ($col in $Cols) {
        $srcCell = $ws.Range("$col$SrcRow")
        $dstCell = $ws.Range("$col$DstRow")
        $dstCell.Value = $srcCell.Value()
    }
}

# Footer row (timestamp / page / developed-by): row 16 -> row 17
Move-RowDown 16 17 @("A", "G", "K")

# Subtotal row: row 15 -> row 16
Move-RowDown 15 16 @("P", "Q")

# Item rows: row 14 -> 15, 13 -> 14, 12 -> 13, 11 -> 12
Move-RowDown 14 15 @("A", "C", "H", "L", "N", "P", "Q")
Move-RowDown 13 14 @("A", "C", "H", "L", "N", "P", "Q")
Move-RowDown 12 13 @("A", "C", "H", "L", "N", "P", "Q")
Move-RowDown 11 12 @("A", "C", "H", "L", "N", "P", "Q")

# ------------------------------------------------------------------
# 2) Re-create the merge pattern one row lower for every moved item row,
#    the subtotal row and the footer row (merges are not part of a
#    format paste, so they have to be restored explicitly).
# ------------------------------------------------------------------

# Un-merge the old locations that have now been vacated / repurposed.
$ws.Range("P15:Q15").UnMerge()
$ws.Range("A16:F16").UnMerge()
$ws.Range("G16:I16").UnMerge()
$ws.Range("K16:Q16").UnMerge()

# Footer row 17
$ws.Range("A17:F17").Merge()
$ws.Range("G17:I17").Merge()
$ws.Range("K17:Q17").Merge()

# Subtotal row 16
$ws.Range("P16:Q16").Merge()

# Item rows 12-15 keep the same merge layout as any other item row.
foreach ($r in 12..15) {
    $ws.Range("A$r`:B$r").Merge()
    $ws.Range("C$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
    $ws.Range("N$r`:O$r").Merge()
}

# ------------------------------------------------------------------
# 3) Populate the freed-up row 11 with the new item, copying the
#    formatting of a normal item row first.
# ------------------------------------------------------------------

$ws.Range("A12:Q12").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "HEMOCLAR 0.5% CREAM 40 GM"
$ws.Range("H11").Value = "0:0"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "46.00"
$ws.Range("P11").Value = "46.0000"
$ws.Range("Q11").Value = "1:0"

# ------------------------------------------------------------------
# 4) Renumber the "م" (item index) column for rows 12-15 (5..8 -> 6..9)
#    and update the subtotal to include the new item's sell price.
# ------------------------------------------------------------------

$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9

$ws.Range("P16").Value = 684.98

# ------------------------------------------------------------------
# 5) Refresh the generated-on timestamp shown in the footer.
# ------------------------------------------------------------------

$ws.Range("A17").Value = "Tuesday, 30 September, 2025 10:32 AM"

Write-Output "Inserted HEMOCLAR row, shifted table, updated subtotal & timestamp."
